$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking Price/Volume columns to be stored as plain text
# (matching the source data export, which always writes these as strings,
# e.g. "-0.52%" / "5.071") rather than letting Excel auto-convert them to
# numbers or percentages on assignment.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '312.98'
$ws.Range("E2").Value = '-0.60%'

$ws.Range("D3").Value = '37.74'
$ws.Range("E3").Value = '-3.88%'

$ws.Range("D4").Value = '5.070'
$ws.Range("E4").Value = '-1.45%'

$ws.Range("D5").Value = '0.07768'
$ws.Range("E5").Value = '-4.91%'

$ws.Range("D6").Value = '4.356'
$ws.Range("E6").Value = '-0.64%'

$ws.Range("D7").Value = '1.903'
$ws.Range("E7").Value = '-4.41%'

$ws.Range("D8").Value = '8.195'
$ws.Range("E8").Value = '-1.81%'

$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '0.9171'
$ws.Range("E9").Value = '-2.08%'

$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '0.1244'
$ws.Range("E10").Value = '-4.61%'

$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1899'
$ws.Range("E11").Value = '-3.47%'

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.08919'
$ws.Range("E12").Value = '-0.68%'

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03388'
$ws.Range("E13").Value = '-3.89%'

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09704'
$ws.Range("E14").Value = '-0.26%'

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001368'
$ws.Range("E15").Value = '-2.83%'

$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '0.005852'
$ws.Range("E16").Value = '-3.79%'

$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.532'
$ws.Range("E17").Value = '-2.58%'

$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '3.012'
$ws.Range("E18").Value = '-3.97%'

$ws.Range("D19").Value = '0.3409'
$ws.Range("E19").Value = '-1.77%'

$ws.Range("E20").Value = '-0.98%'

$ws.Range("D21").Value = '5.031'
$ws.Range("E21").Value = '0.80%'

$ws.Range("D22").Value = '0.2592'
$ws.Range("E22").Value = '4.14%'

$ws.Range("E23").Value = '5,592.31%'

$ws.Range("E24").Value = '0.93%'

$ws.Range("E25").Value = '-2.25%'

$ws.Range("D26").Value = '0.004241'
$ws.Range("E26").Value = '-10.94%'

$ws.Range("E27").Value = '-65.28%'

$ws.Range("D39").Value = '0.02136'
$ws.Range("E39").Value = '-4.18%'

$ws.Range("D40").Value = '0.04970'
$ws.Range("E40").Value = '-4.19%'

$ws.Range("D41").Value = '0.007803'
$ws.Range("E41").Value = '0.86%'

$ws.Range("D42").Value = '0.009876'
$ws.Range("E42").Value = '-4.13%'

$ws.Range("E43").Value = '-3.83%'

$ws.Range("D44").Value = '0.002062'
$ws.Range("E44").Value = '-1.86%'

$ws.Range("D45").Value = '0.009688'
$ws.Range("E45").Value = '15.20%'

$ws.Range("D46").Value = '0.00006518'
$ws.Range("E46").Value = '-4.41%'

$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").Value = '0.04%'

$ws.Range("D48").Value = '0.003056'
$ws.Range("E48").Value = '1.62%'

$ws.Range("E49").Value = '-0.11%'

$ws.Range("D50").Value = '0.00002102'
$ws.Range("E50").Value = '0.04%'

$ws.Range("E51").Value = '0.04%'
